$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: update the reporting period to 2022 H1, with new validation/update dates ---
$ws.Range("A8").Value = 2022
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "06/30/2022"
$ws.Range("I8").Value = "07/11/2022"
$ws.Range("J8").Value = "07/11/2022"

# H8 keeps the same text ("Secretaria Administrativa (UPP)") but picks up a fresh
# paste-in style: Calibri 11 black text, bordered, no special alignment.
$ws.Range("H8").ClearFormats()
$ws.Range("H8").Borders.LineStyle = 1
$ws.Range("H8").Font.Name = "Calibri"
$ws.Range("H8").Font.Size = 11
$ws.Range("H8").Font.Color = 0

# K8 gets the updated note text
$ws.Range("K8").Value = "Los criterios e hipervínculos que se encuentran en blanco es porque para este periodo no se realizo ninguna baja de bienes inmuebles."

# Row height shrinks back down now that the note text is shorter
$ws.Rows(8).RowHeight = 60

# Column K is very slightly narrower
$ws.Columns(11).ColumnWidth = 38.140625

# Selection moved to C18 before saving
$ws.Range("C18").Select()
